# Fruta / hortaliza, semanal
# Insert two new weekly price rows for "Choclo" (row 151 and 152), shifting
# the existing rows 151-193 down to 153-195.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 151 (this pushes old rows
# 151-193 down to become rows 153-195, matching the new dimension A1:R195).
$ws.Range("A151:A152").EntireRow.Insert()

# Populate the first new row (151) - Choclero, Primera, Región del Maule
$ws.Range("A151").Value = 7
$ws.Range("B151").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C151").Value = "Ñuble"
$ws.Range("D151").Value = 44642
$ws.Range("E151").Value = 16
$ws.Range("F151").Value = 100112024
$ws.Range("G151").Value = "Choclo"
$ws.Range("H151").Value = "Choclero"
$ws.Range("I151").Value = "Primera"
$ws.Range("J151").Value = 10000
$ws.Range("K151").Value = 200
$ws.Range("L151").Value = 200
$ws.Range("M151").Value = 200
$ws.Range("N151").Value = "$/unidad"
$ws.Range("O151").Value = "Región del Maule"
$ws.Range("P151").Value = 200
$ws.Range("Q151").Value = 1
$ws.Range("R151").Value = "Hortaliza"

# Populate the second new row (152) - Choclero, Segunda, Región del Maule
$ws.Range("A152").Value = 7
$ws.Range("B152").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C152").Value = "Ñuble"
$ws.Range("D152").Value = 44642
$ws.Range("E152").Value = 16
$ws.Range("F152").Value = 100112024
$ws.Range("G152").Value = "Choclo"
$ws.Range("H152").Value = "Choclero"
$ws.Range("I152").Value = "Segunda"
$ws.Range("J152").Value = 10000
$ws.Range("K152").Value = 150
$ws.Range("L152").Value = 150
$ws.Range("M152").Value = 150
$ws.Range("N152").Value = "$/unidad"
$ws.Range("O152").Value = "Región del Maule"
$ws.Range("P152").Value = 150
$ws.Range("Q152").Value = 1
$ws.Range("R152").Value = "Hortaliza"
